$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 20.69072717484972
$ws.Cells.Item(2, 3).Value = 7.952541620753292
$ws.Cells.Item(2, 4).Value = 13.2751333515032
$ws.Cells.Item(2, 5).Value = 13.38522071206944
$ws.Cells.Item(2, 7).Value = 3.697049035650651
$ws.Cells.Item(2, 9).Value = 28.17620143573518
$ws.Cells.Item(2, 10).Value = 8.187482756625961
$ws.Cells.Item(2, 12).Value = 13.0202266327786
$ws.Cells.Item(2, 15).Value = 30.25265188611007
$ws.Cells.Item(3, 2).Value = 20.18841411796762
$ws.Cells.Item(3, 3).Value = 7.515664839384573
$ws.Cells.Item(3, 4).Value = 13.28769929825645
$ws.Cells.Item(3, 5).Value = 13.42499171706828
$ws.Cells.Item(3, 7).Value = 3.69965029173114
$ws.Cells.Item(3, 9).Value = 28.3336482195954
$ws.Cells.Item(3, 10).Value = 8.195660685478027
$ws.Cells.Item(3, 12).Value = 13.00194378577044
$ws.Cells.Item(3, 15).Value = 30.3723061299843
$ws.Cells.Item(4, 2).Value = 19.87577909099843
$ws.Cells.Item(4, 3).Value = 7.232781118002418
$ws.Cells.Item(4, 4).Value = 13.29788655317368
$ws.Cells.Item(4, 5).Value = 13.45115326932511
$ws.Cells.Item(4, 7).Value = 3.701331861627673
$ws.Cells.Item(4, 9).Value = 28.4367059193324
$ws.Cells.Item(4, 10).Value = 8.200970172940979
$ws.Cells.Item(4, 12).Value = 12.99252058032478
$ws.Cells.Item(4, 15).Value = 30.45334626298173
$ws.Cells.Item(5, 2).Value = 19.74749282855333
$ws.Cells.Item(5, 3).Value = 7.113879922479143
$ws.Cells.Item(5, 4).Value = 13.30265938585486
$ws.Cells.Item(5, 5).Value = 13.46225280425366
$ws.Cells.Item(5, 7).Value = 3.702038405824173
$ws.Cells.Item(5, 9).Value = 28.48030676829028
$ws.Cells.Item(5, 10).Value = 8.203206521091111
$ws.Cells.Item(5, 12).Value = 12.98913660559778
$ws.Cells.Item(5, 15).Value = 30.48826856540904
$ws.Cells.Item(6, 2).Value = 19.72614280141868
$ws.Cells.Item(6, 3).Value = 7.093919539464352
$ws.Cells.Item(6, 4).Value = 13.30348943858475
$ws.Cells.Item(6, 5).Value = 13.46412237064967
$ws.Cells.Item(6, 7).Value = 3.702157014941517
$ws.Cells.Item(6, 9).Value = 28.48764349180543
$ws.Cells.Item(6, 10).Value = 8.203582261619784
$ws.Cells.Item(6, 12).Value = 12.98860232071052
$ws.Cells.Item(6, 15).Value = 30.49418182434639
$ws.Cells.Item(7, 2).Value = 19.8740523242518
$ws.Cells.Item(7, 3).Value = 7.231192162011192
$ws.Cells.Item(7, 4).Value = 13.29794840531265
$ws.Cells.Item(7, 5).Value = 13.45130118550357
$ws.Cells.Item(7, 7).Value = 3.701341304031134
$ws.Cells.Item(7, 9).Value = 28.43728744305488
$ws.Cells.Item(7, 10).Value = 8.201000038511912
$ws.Cells.Item(7, 12).Value = 12.99247309279145
$ws.Cells.Item(7, 15).Value = 30.45380956167269
$ws.Cells.Item(8, 2).Value = 20.51850598521946
$ws.Cells.Item(8, 3).Value = 7.804975582878051
$ws.Cells.Item(8, 4).Value = 13.27895308350092
$ws.Cells.Item(8, 5).Value = 13.39857247521608
$ws.Cells.Item(8, 7).Value = 3.697928476635498
$ws.Cells.Item(8, 9).Value = 28.22916301905131
$ws.Cells.Item(8, 10).Value = 8.190242820074666
$ws.Cells.Item(8, 12).Value = 13.01355033839617
$ws.Cells.Item(8, 15).Value = 30.29233339998792
$ws.Cells.Item(9, 2).Value = 21.74163793717575
$ws.Cells.Item(9, 3).Value = 8.812073562699023
$ws.Cells.Item(9, 4).Value = 13.2613132107176
$ws.Cells.Item(9, 5).Value = 13.3089732918248
$ws.Cells.Item(9, 7).Value = 3.691902287146123
$ws.Cells.Item(9, 9).Value = 27.8717598245363
$ws.Cells.Item(9, 10).Value = 8.17142488858393
$ws.Cells.Item(9, 12).Value = 13.0690532929507
$ws.Cells.Item(9, 15).Value = 30.0360331492904
$ws.Cells.Item(10, 2).Value = 22.60652010395209
$ws.Cells.Item(10, 3).Value = 9.478085960658442
$ws.Cells.Item(10, 4).Value = 13.26029621329069
$ws.Cells.Item(10, 5).Value = 13.25153068462993
$ws.Cells.Item(10, 7).Value = 3.687876555275492
$ws.Cells.Item(10, 9).Value = 27.64019659713573
$ws.Cells.Item(10, 10).Value = 8.158973748580161
$ws.Cells.Item(10, 12).Value = 13.11828041595935
$ws.Cells.Item(10, 15).Value = 29.88489022964821
$ws.Cells.Item(11, 2).Value = 22.990953996889
$ws.Cells.Item(11, 3).Value = 9.76477257957832
$ws.Cells.Item(11, 4).Value = 13.26242038133705
$ws.Cells.Item(11, 5).Value = 13.22721367141487
$ws.Cells.Item(11, 7).Value = 3.686131410761507
$ws.Cells.Item(11, 9).Value = 27.54161126995355
$ws.Cells.Item(11, 10).Value = 8.15360493574422
$ws.Cells.Item(11, 12).Value = 13.14246430515301
$ws.Cells.Item(11, 15).Value = 29.82427635168045
$ws.Cells.Item(12, 2).Value = 23.13510711427503
$ws.Cells.Item(12, 3).Value = 9.870977363181991
$ws.Cells.Item(12, 4).Value = 13.26359578641327
$ws.Cells.Item(12, 5).Value = 13.21826587526146
$ws.Cells.Item(12, 7).Value = 3.685482889463597
$ws.Cells.Item(12, 9).Value = 27.50525308709377
$ws.Cells.Item(12, 10).Value = 8.151614144654932
$ws.Cells.Item(12, 12).Value = 13.15187525630136
$ws.Cells.Item(12, 15).Value = 29.80250024539144
$ws.Cells.Item(13, 2).Value = 23.10412641572924
$ws.Cells.Item(13, 3).Value = 9.848209297326713
$ws.Cells.Item(13, 4).Value = 13.26332615773864
$ws.Cells.Item(13, 5).Value = 13.22018136102558
$ws.Cells.Item(13, 7).Value = 3.685622012920708
$ws.Cells.Item(13, 9).Value = 27.51304011114344
$ws.Cells.Item(13, 10).Value = 8.152041020666932
$ws.Cells.Item(13, 12).Value = 13.14983725820768
$ws.Cells.Item(13, 15).Value = 29.80713767190754
$ws.Cells.Item(14, 2).Value = 23.00284279244207
$ws.Cells.Item(14, 3).Value = 9.773557425060098
$ws.Cells.Item(14, 4).Value = 13.26250965286559
$ws.Cells.Item(14, 5).Value = 13.22647231120188
$ws.Cells.Item(14, 7).Value = 3.686077809846198
$ws.Cells.Item(14, 9).Value = 27.53860052504043
$ws.Cells.Item(14, 10).Value = 8.153440306302231
$ws.Cells.Item(14, 12).Value = 13.14323350454749
$ws.Cells.Item(14, 15).Value = 29.82246119316438
$ws.Cells.Item(15, 2).Value = 22.94061460743494
$ws.Cells.Item(15, 3).Value = 9.72752364058616
$ws.Cells.Item(15, 4).Value = 13.26205780704379
$ws.Cells.Item(15, 5).Value = 13.23035961931734
$ws.Cells.Item(15, 7).Value = 3.686358602028406
$ws.Cells.Item(15, 9).Value = 27.55438393570821
$ws.Cells.Item(15, 10).Value = 8.154302907042007
$ws.Cells.Item(15, 12).Value = 13.13922133314607
$ws.Cells.Item(15, 15).Value = 29.83200076238876
$ws.Cells.Item(16, 2).Value = 22.58120399156579
$ws.Cells.Item(16, 3).Value = 9.459021187032841
$ws.Cells.Item(16, 4).Value = 13.26020937071383
$ws.Cells.Item(16, 5).Value = 13.25315633261382
$ws.Cells.Item(16, 7).Value = 3.687992333132529
$ws.Cells.Item(16, 9).Value = 27.64677551738309
$ws.Cells.Item(16, 10).Value = 8.159330537802905
$ws.Cells.Item(16, 12).Value = 13.1167355956169
$ws.Cells.Item(16, 15).Value = 29.88901586355352
$ws.Cells.Item(17, 2).Value = 22.3583158053709
$ws.Cells.Item(17, 3).Value = 9.290120033834144
$ws.Cells.Item(17, 4).Value = 13.25973734988779
$ws.Cells.Item(17, 5).Value = 13.26760574263924
$ws.Cells.Item(17, 7).Value = 3.689016600421172
$ws.Cells.Item(17, 9).Value = 27.70518651120603
$ws.Cells.Item(17, 10).Value = 8.16249031148083
$ws.Cells.Item(17, 12).Value = 13.10339679302097
$ws.Cells.Item(17, 15).Value = 29.92608261505449
$ws.Cells.Item(18, 2).Value = 22.22927580845086
$ws.Cells.Item(18, 3).Value = 9.191440254241281
$ws.Cells.Item(18, 4).Value = 13.25970940890516
$ws.Cells.Item(18, 5).Value = 13.27608741364135
$ws.Cells.Item(18, 7).Value = 3.689613847545074
$ws.Cells.Item(18, 9).Value = 27.73941851185299
$ws.Cells.Item(18, 10).Value = 8.164335533971453
$ws.Cells.Item(18, 12).Value = 13.09589344940658
$ws.Cells.Item(18, 15).Value = 29.94816841077462
$ws.Cells.Item(19, 2).Value = 22.18544480241738
$ws.Cells.Item(19, 3).Value = 9.157766346115178
$ws.Cells.Item(19, 4).Value = 13.25974180444464
$ws.Cells.Item(19, 5).Value = 13.27898849818522
$ws.Cells.Item(19, 7).Value = 3.689817461042138
$ws.Cells.Item(19, 9).Value = 27.7511179529433
$ws.Cells.Item(19, 10).Value = 8.164965076072876
$ws.Cells.Item(19, 12).Value = 13.09338206476135
$ws.Cells.Item(19, 15).Value = 29.95577764539009
$ws.Cells.Item(20, 2).Value = 22.38213059876524
$ws.Cells.Item(20, 3).Value = 9.308258563276391
$ws.Cells.Item(20, 4).Value = 13.25976239670348
$ws.Cells.Item(20, 5).Value = 13.26604990829545
$ws.Cells.Item(20, 7).Value = 3.688906725955249
$ws.Cells.Item(20, 9).Value = 27.69890276555777
$ws.Cells.Item(20, 10).Value = 8.162151071982956
$ws.Cells.Item(20, 12).Value = 13.10479929401131
$ws.Cells.Item(20, 15).Value = 29.92205747384078
$ws.Cells.Item(21, 2).Value = 23.03263185868276
$ws.Cells.Item(21, 3).Value = 9.79554856366506
$ws.Cells.Item(21, 4).Value = 13.26273941917542
$ws.Cells.Item(21, 5).Value = 13.22461743839427
$ws.Cells.Item(21, 7).Value = 3.685943597209016
$ws.Cells.Item(21, 9).Value = 27.53106635699869
$ws.Cells.Item(21, 10).Value = 8.153028156788988
$ws.Cells.Item(21, 12).Value = 13.14516635592709
$ws.Cells.Item(21, 15).Value = 29.81792830959403
$ws.Cells.Item(22, 2).Value = 23.44942158536363
$ws.Cells.Item(22, 3).Value = 10.1002769843646
$ws.Cells.Item(22, 4).Value = 13.26684712091626
$ws.Cells.Item(22, 5).Value = 13.19905734768358
$ws.Cells.Item(22, 7).Value = 3.684078842628144
$ws.Cells.Item(22, 9).Value = 27.42705379385362
$ws.Cells.Item(22, 10).Value = 8.147312047788912
$ws.Cells.Item(22, 12).Value = 13.17302125341205
$ws.Cells.Item(22, 15).Value = 29.75673742066623
$ws.Cells.Item(23, 2).Value = 23.22777693373491
$ws.Cells.Item(23, 3).Value = 9.93889918587786
$ws.Cells.Item(23, 4).Value = 13.26445730590213
$ws.Cells.Item(23, 5).Value = 13.21256041390826
$ws.Cells.Item(23, 7).Value = 3.68506754717164
$ws.Cells.Item(23, 9).Value = 27.48204676193708
$ws.Cells.Item(23, 10).Value = 8.150340376466382
$ws.Cells.Item(23, 12).Value = 13.15802133241976
$ws.Cells.Item(23, 15).Value = 29.7887661141724
$ws.Cells.Item(24, 2).Value = 22.37136671973434
$ws.Cells.Item(24, 3).Value = 9.300063042682453
$ws.Cells.Item(24, 4).Value = 13.25975031471495
$ws.Cells.Item(24, 5).Value = 13.26675275738318
$ws.Cells.Item(24, 7).Value = 3.688956374085107
$ws.Cells.Item(24, 9).Value = 27.701741620057
$ws.Cells.Item(24, 10).Value = 8.162304352973711
$ws.Cells.Item(24, 12).Value = 13.10416470805541
$ws.Cells.Item(24, 15).Value = 29.92387482454577
$ws.Cells.Item(25, 2).Value = 21.4160654809149
$ws.Cells.Item(25, 3).Value = 8.552498046916746
$ws.Cells.Item(25, 4).Value = 13.26398629721947
$ws.Cells.Item(25, 5).Value = 13.33173759926146
$ws.Cells.Item(25, 7).Value = 3.693461661921299
$ws.Cells.Item(25, 9).Value = 27.96300662918972
$ws.Cells.Item(25, 10).Value = 8.176273296881536
$ws.Cells.Item(25, 12).Value = 13.05253978295331
$ws.Cells.Item(25, 15).Value = 30.09886932192485
